$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header cell + column data for the "Presskit" column (I)
$ws.Range("I1").Value = "Presskit"
$ws.Range("I2").Value = "not important enough"
$ws.Range("I3").Value = "not important enough"
$ws.Range("I4").Value = "not important enough"
$ws.Range("I5").Value = "small link"
$ws.Range("I6").Value = "not important enough"
$ws.Range("I7").Value = "not important enough"
$ws.Range("I9").Value = "not important enough"

# New row 10: Presskit
$ws.Range("A10").Value = "Presskit"
$ws.Range("I10").Value = "x"

# New header cell for the "Twitch page" column (J)
$ws.Range("J1").Value = "Twitch page"

# New row 11: Twitch page
$ws.Range("A11").Value = "Twitch page"

# Update selection to match the saved state
$ws.Range("J5").Select()
